$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 104.23077
$ws.Range("I9").Value = 102.333336
$ws.Range("J9").Value = 108.5
$ws.Range("K9").Value = 102.333336
$ws.Range("L9").Value = 108.5
$ws.Range("M9").Value = 66.666664
$ws.Range("N9").Value = -446.5
$ws.Range("H17").Value = 4333.3335
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 9000
$ws.Range("M17").Value = -8832
$ws.Range("H29").Value = 2314.6667
$ws.Range("I29").Value = 195.5
$ws.Range("J29").Value = 3374.25
$ws.Range("K29").Value = 586.5
$ws.Range("L29").Value = 10122.75
$ws.Range("M29").Value = -305.5
$ws.Range("N29").Value = -10684.75
$ws.Range("H38").Value = 174.57143
$ws.Range("I38").Value = 174.57143
$ws.Range("K38").Value = 523.71429
$ws.Range("M38").Value = -151.71429
$ws.Range("H40").Value = 2927.625
$ws.Range("I40").Value = 2330
$ws.Range("K40").Value = 2330
$ws.Range("M40").Value = -2155
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H100").Value = 30000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 30000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 30000
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -31082
$ws.Range("H101").Value = 997.1667
$ws.Range("I101").Value = 596.6
$ws.Range("K101").Value = 1789.8
$ws.Range("M101").Value = -167.8000000000002
$ws.Range("H106").Value = 8600
$ws.Range("I106").Value = 8600
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 8600
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -7969
$ws.Range("N106").Value = ""
$ws.Range("H111").Value = 768.5
$ws.Range("I111").Value = 768.5
$ws.Range("K111").Value = 2305.5
$ws.Range("M111").Value = 761.5
$ws.Range("H132").Value = 83336360
$ws.Range("I132").Value = 90911850
$ws.Range("K132").Value = 272735550
$ws.Range("M132").Value = -272733020
$ws.Range("H135").Value = 1627.2727
$ws.Range("I135").Value = 1237.75
$ws.Range("J135").Value = 2666
$ws.Range("K135").Value = 11139.75
$ws.Range("L135").Value = 23994
$ws.Range("M135").Value = -8604.75
$ws.Range("N135").Value = -29064
$ws.Range("H141").Value = 5939.923
$ws.Range("I141").Value = 5601.8335
$ws.Range("K141").Value = 16805.5005
$ws.Range("M141").Value = -11625.5005

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 413.75
$ws.Range("I2").Value = 413.75
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 413.75
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -300.75
$ws.Range("N2").Value = ""
$ws.Range("H30").Value = 1759.6
$ws.Range("I30").Value = 933.3333
$ws.Range("J30").Value = 2999
$ws.Range("K30").Value = 933.3333
$ws.Range("L30").Value = 2999
$ws.Range("M30").Value = -783.3333
$ws.Range("N30").Value = -3299
$ws.Range("H35").Value = 3374.5
$ws.Range("J35").Value = 4749.5
$ws.Range("L35").Value = 4749.5
$ws.Range("N35").Value = -5561.5
$ws.Range("H39").Value = 1499
$ws.Range("I39").Value = 1499
$ws.Range("K39").Value = 1499
$ws.Range("M39").Value = -979
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50972
$ws.Range("H45").Value = 1650
$ws.Range("I45").Value = 1650
$ws.Range("K45").Value = 1650
$ws.Range("M45").Value = -1273
$ws.Range("H74").Value = 3550
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4126
$ws.Range("H77").Value = 3550
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20632
$ws.Range("H102").Value = 2499.8
$ws.Range("I102").Value = 2499.75
$ws.Range("K102").Value = 2499.75
$ws.Range("M102").Value = -877.75
$ws.Range("H116").Value = 413.75
$ws.Range("I116").Value = 413.75
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 413.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1880.25
$ws.Range("N116").Value = ""
$ws.Range("H132").Value = 3763.5
$ws.Range("I132").Value = 3763.5
$ws.Range("K132").Value = 11290.5
$ws.Range("M132").Value = -8760.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 413.75
$ws.Range("I3").Value = 413.75
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 413.75
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -299.75
$ws.Range("N3").Value = ""
$ws.Range("H134").Value = 6037.0527
$ws.Range("I134").Value = 6150.222
$ws.Range("K134").Value = 18450.666
$ws.Range("M134").Value = -15915.666

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 214.76923
$ws.Range("I5").Value = 157.45454
$ws.Range("J5").Value = 530
$ws.Range("K5").Value = 157.45454
$ws.Range("L5").Value = 530
$ws.Range("M5").Value = -45.45454000000001
$ws.Range("N5").Value = -754
$ws.Range("H35").Value = 194.15384
$ws.Range("I35").Value = 194.15384
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 194.15384
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 99.84616
$ws.Range("N35").Value = ""
$ws.Range("H105").Value = 1522.8
$ws.Range("I105").Value = 1522.8
$ws.Range("K105").Value = 1522.8
$ws.Range("M105").Value = 224.2
$ws.Range("H137").Value = 60000
$ws.Range("I137").Value = 60000
$ws.Range("K137").Value = 60000
$ws.Range("M137").Value = -54900

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 786.6667
$ws.Range("I18").Value = 786.6667
$ws.Range("K18").Value = 2360.0001
$ws.Range("M18").Value = -2191.0001
$ws.Range("H70").Value = 3210.6
$ws.Range("I70").Value = 2017.6666
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 6052.9998
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -5737.9998
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 3210.6
$ws.Range("I73").Value = 2017.6666
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 6052.9998
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -4960.9998
$ws.Range("N73").Value = -17184
$ws.Range("H80").Value = 1664.6666
$ws.Range("I80").Value = 1999.5
$ws.Range("K80").Value = 5998.5
$ws.Range("M80").Value = -5062.5
$ws.Range("H83").Value = 1664.6666
$ws.Range("I83").Value = 1999.5
$ws.Range("K83").Value = 17995.5
$ws.Range("M83").Value = -13315.5
$ws.Range("H107").Value = 1861.875
$ws.Range("J107").Value = 1919.4
$ws.Range("L107").Value = 5758.200000000001
$ws.Range("N107").Value = -9598.200000000001
$ws.Range("H120").Value = 24444.223
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 24444.223
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 73332.66900000001
$ws.Range("M120").Value = ""
$ws.Range("N120").Value = -83008.66900000001
$ws.Range("H121").Value = 1009.1429
$ws.Range("J121").Value = 1009.1429
$ws.Range("L121").Value = 3027.4287
$ws.Range("N121").Value = -5647.4287
$ws.Range("H131").Value = 2319.5334
$ws.Range("I131").Value = 1370.4286
$ws.Range("K131").Value = 4111.2858
$ws.Range("M131").Value = 928.7142000000003

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""
$ws.Range("H27").Value = 12999.667
$ws.Range("J27").Value = 14499.5
$ws.Range("L27").Value = 14499.5
$ws.Range("N27").Value = -14831.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10272

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -31240
$ws.Range("H68").Value = 47500
$ws.Range("J68").Value = 47500
$ws.Range("L68").Value = 47500
$ws.Range("N68").Value = -49122
$ws.Range("H71").Value = 47500
$ws.Range("J71").Value = 47500
$ws.Range("L71").Value = 142500
$ws.Range("N71").Value = -150612
$ws.Range("H100").Value = 793
$ws.Range("I100").Value = 723.3333
$ws.Range("J100").Value = 932.3333
$ws.Range("K100").Value = 1446.6666
$ws.Range("L100").Value = 1864.6666
$ws.Range("M100").Value = -905.6666
$ws.Range("N100").Value = -2946.6666
$ws.Range("H113").Value = 1632.1538
$ws.Range("I113").Value = 992.875
$ws.Range("K113").Value = 2978.625
$ws.Range("M113").Value = -808.625
$ws.Range("H122").Value = 852.8570999999999
$ws.Range("I122").Value = 852.8570999999999
$ws.Range("K122").Value = 2558.5713
$ws.Range("M122").Value = -108.5712999999996
$ws.Range("H136").Value = 2151.95
$ws.Range("I136").Value = 2081.0527
$ws.Range("J136").Value = 3499
$ws.Range("K136").Value = 6243.158100000001
$ws.Range("L136").Value = 10497
$ws.Range("M136").Value = -3693.158100000001
$ws.Range("N136").Value = -15597

Write-Host "Updated cells: 253 set, 9 cleared"